$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 98524.67999999999
$ws.Range("I15").Value = 98524.67999999999
$ws.Range("K15").Value = 295574.04
$ws.Range("M15").Value = -295405.04

$ws.Range("H33").Value = 185.94444
$ws.Range("I33").Value = 191.46666
$ws.Range("K33").Value = 191.46666
$ws.Range("M33").Value = 37.53334000000001

$ws.Range("H70").Value = 2058.611
$ws.Range("I70").Value = 1536.2727
$ws.Range("J70").Value = 2879.4285
$ws.Range("K70").Value = 4608.8181
$ws.Range("L70").Value = 8638.2855
$ws.Range("M70").Value = -4338.8181
$ws.Range("N70").Value = -9178.2855

$ws.Range("H73").Value = 2058.611
$ws.Range("I73").Value = 1536.2727
$ws.Range("J73").Value = 2879.4285
$ws.Range("K73").Value = 4608.8181
$ws.Range("L73").Value = 8638.2855
$ws.Range("M73").Value = -3672.8181
$ws.Range("N73").Value = -10510.2855

$ws.Range("H129").Value = 1388.6364
$ws.Range("J129").Value = 1928.4286
$ws.Range("L129").Value = 5785.2858
$ws.Range("N129").Value = -15785.2858

$ws.Range("H137").Value = 55557280
$ws.Range("I137").Value = 76924460
$ws.Range("J137").Value = 2600
$ws.Range("K137").Value = 230773380
$ws.Range("L137").Value = 7800
$ws.Range("M137").Value = -230770830
$ws.Range("N137").Value = -12900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23380.613
$ws.Range("I32").Value = 5594.2
$ws.Range("J32").Value = 512507
$ws.Range("K32").Value = 5594.2
$ws.Range("L32").Value = 512507
$ws.Range("M32").Value = -5307.2
$ws.Range("N32").Value = -513081

$ws.Range("H61").Value = 1892.3667
$ws.Range("I61").Value = 1623.7916
$ws.Range("K61").Value = 1623.7916
$ws.Range("M61").Value = -1411.7916

$ws.Range("H132").Value = 2179.5686
$ws.Range("I132").Value = 2189.2778
$ws.Range("J132").Value = 2156.2666
$ws.Range("K132").Value = 6567.8334
$ws.Range("L132").Value = 6468.7998
$ws.Range("M132").Value = -4037.8334
$ws.Range("N132").Value = -11528.7998

$ws.Range("H136").Value = 1892.3667
$ws.Range("I136").Value = 1623.7916
$ws.Range("K136").Value = 4871.3748
$ws.Range("M136").Value = -2321.3748

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 5397.4546
$ws.Range("I19").Value = 337.14285
$ws.Range("J19").Value = 14253
$ws.Range("K19").Value = 337.14285
$ws.Range("L19").Value = 14253
$ws.Range("M19").Value = -167.14285
$ws.Range("N19").Value = -14593

$ws.Range("H24").Value = 5397.4546
$ws.Range("I24").Value = 337.14285
$ws.Range("J24").Value = 14253
$ws.Range("K24").Value = 337.14285
$ws.Range("L24").Value = 14253
$ws.Range("M24").Value = -167.14285
$ws.Range("N24").Value = -14593

$ws.Range("H28").Value = 30000
$ws.Range("J28").Value = 30000
$ws.Range("L28").Value = 30000
$ws.Range("N28").Value = -30490

$ws.Range("H31").Value = 3515.625
$ws.Range("I31").Value = 1616.5358
$ws.Range("K31").Value = 1616.5358
$ws.Range("M31").Value = -1321.5358

$ws.Range("H34").Value = 3515.625
$ws.Range("I34").Value = 1616.5358
$ws.Range("K34").Value = 1616.5358
$ws.Range("M34").Value = -1414.5358

$ws.Range("H74").Value = 18546
$ws.Range("J74").Value = 18546
$ws.Range("L74").Value = 18546
$ws.Range("N74").Value = -20294

$ws.Range("H77").Value = 18546
$ws.Range("J77").Value = 18546
$ws.Range("L77").Value = 55638
$ws.Range("N77").Value = -64374

$ws.Range("H99").Value = 20835920
$ws.Range("I99").Value = 2110
$ws.Range("J99").Value = 41669732
$ws.Range("K99").Value = 2110
$ws.Range("L99").Value = 41669732
$ws.Range("M99").Value = -612
$ws.Range("N99").Value = -41672728

$ws.Range("H126").Value = 20835920
$ws.Range("I126").Value = 2110
$ws.Range("J126").Value = 41669732
$ws.Range("K126").Value = 6330
$ws.Range("L126").Value = 125009196
$ws.Range("M126").Value = -3860
$ws.Range("N126").Value = -125014136

$ws.Range("H134").Value = 30614500
$ws.Range("I134").Value = 43479936
$ws.Range("J134").Value = 19233538
$ws.Range("K134").Value = 130439808
$ws.Range("L134").Value = 57700614
$ws.Range("M134").Value = -130437273
$ws.Range("N134").Value = -57705684

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1310
$ws.Range("J80").Value = 1310
$ws.Range("L80").Value = 3930
$ws.Range("N80").Value = -5802

$ws.Range("H83").Value = 1310
$ws.Range("J83").Value = 1310
$ws.Range("L83").Value = 11790
$ws.Range("N83").Value = -21150

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = -888
$ws.Range("N5").Value = -1224

$ws.Range("H6").Value = 29909
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 29909
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 29909
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -30135

$ws.Range("H16").Value = 29909
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 29909
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 29909
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -30409

$ws.Range("H58").Value = 1000
$ws.Range("I58").Value = 1000
$ws.Range("K58").Value = 1000
$ws.Range("M58").Value = -723

$ws.Range("H122").Value = 1939.3125
$ws.Range("I122").Value = 1547.1538
$ws.Range("J122").Value = 3638.6667
$ws.Range("K122").Value = 4641.4614
$ws.Range("L122").Value = 10916.0001
$ws.Range("M122").Value = -2191.4614
$ws.Range("N122").Value = -15816.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3072.8696
$ws.Range("I7").Value = 2213.8572
$ws.Range("J7").Value = 3448.6875
$ws.Range("K7").Value = 2213.8572
$ws.Range("L7").Value = 3448.6875
$ws.Range("M7").Value = -2101.8572
$ws.Range("N7").Value = -3672.6875

$ws.Range("H126").Value = 3072.8696
$ws.Range("I126").Value = 2213.8572
$ws.Range("J126").Value = 3448.6875
$ws.Range("K126").Value = 6641.571599999999
$ws.Range("L126").Value = 10346.0625
$ws.Range("M126").Value = -4171.571599999999
$ws.Range("N126").Value = -15286.0625

$ws.Range("H138").Value = 70429
$ws.Range("J138").Value = 70429
$ws.Range("L138").Value = 70429
$ws.Range("N138").Value = -80709

$ws.Range("H141").Value = 58383
$ws.Range("J141").Value = 58383
$ws.Range("L141").Value = 58383
$ws.Range("N141").Value = -68743

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 11309.8
$ws.Range("I9").Value = 16916.666
$ws.Range("J9").Value = 2899.5
$ws.Range("K9").Value = 16916.666
$ws.Range("L9").Value = 2899.5
$ws.Range("M9").Value = -16776.666
$ws.Range("N9").Value = -3179.5

$ws.Range("H103").Value = 517801
$ws.Range("J103").Value = 517801
$ws.Range("L103").Value = 517801
$ws.Range("N103").Value = -520145

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H122").Value = 2074.125
$ws.Range("I122").Value = 1816.6364
$ws.Range("J122").Value = 2640.6
$ws.Range("K122").Value = 5449.9092
$ws.Range("L122").Value = 7921.799999999999
$ws.Range("M122").Value = -2999.9092
$ws.Range("N122").Value = -12821.8

$ws.Range("H137").Value = 61569.168
$ws.Range("J137").Value = 98138.336
$ws.Range("L137").Value = 98138.336
$ws.Range("N137").Value = -108338.336
